$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column B - shifts B:U to C:V
$ws.Columns("B:B").Insert()

# Populate the new column B with header "Level" and row values "Object"/"Component"
$ws.Range("B1").Value = "Level"
$ws.Range("B2").Value = "Object"
$ws.Range("B3").Value = "Component"
$ws.Range("B4").Value = "Component"

$ws.Range("B4").Select()
